# Updated symbol list: refresh the "Price" column (column D) on Sheet1
# with new coin prices. Values are stored as text (matching the
# workbook's existing inline-string convention for column D), so each
# cell is explicitly formatted as Text before the value is written --
# otherwise Excel would auto-convert the numeric-looking string to a
# real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "271.69"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.04"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.378"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06303"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.657"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.747"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.396"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8348"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1625"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08386"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03478"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03125"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09320"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.946"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001720"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04859"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006236"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005472"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001090"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.741"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.320"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01385"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3380"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002684"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04685"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006903"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1175"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003461"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01260"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006261"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7894"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1143"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01241"
